# toilet_coop_2022-09-08.xlsx — refresh scrape snapshot
#
# The crawler re-ran later the same day (new timestamp 2022-09-08 21:00:24
# instead of 07:11:03) and the product list shifted: the two rows for
# "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen" (row 3) and
# "Prix Garantie feuchtes Toilettenpapier 2x70 Stück" (row 4) are no longer
# present, so every following row moves up by two and the timestamp column
# is refreshed for all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two rows that disappeared from the newer crawl; remaining rows
# shift up automatically.
$ws.Rows("3:4").Delete()

# Refresh the timestamp column (O) for the title row (row 2) and every
# product row (rows 3-31) to the new crawl time.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-09-08 21:00:24"
}
